$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from SCD0338 to SCD0026
$ws.Name = "SCD0026"

# Update the test case ID text from SCD0338-011 to SCD0026-011
$ws.Range("B2").Value = "SCD0026-011"
$ws.Range("B3").Value = "SCD0026-011"

# Move the active selection to B4 (and scroll the view back to show column A)
[void]$ws.Range("B4").Select()
